$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.331.85'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.587.40'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '210.06'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').Value = '19.48'
$ws.Range('E10').Value = '  -0.25%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '1.811.62'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.07'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.576.45'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.30'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').Value = '26.332.50'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('D20').Value = '210.85'
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  -3.58%  '
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').Value = '144.97'
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').Value = '7.05'
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  -0.42%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').Value = '1.301.89'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.611'
$ws.Range('E35').Value = '  +2.37%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.44'
$ws.Range('E36').Value = '  -1.72%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('E39').Value = '  -13.00%  '
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('D42').Value = '5.61'
$ws.Range('E42').Value = '  +3.50%  '
$ws.Range('D43').Value = '0.767'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '62.47'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.13'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('D46').Value = '1.723.46'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '87.77'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('E48').Value = '  -5.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0980'
$ws.Range('E50').Value = '  -4.34%  '
$ws.Range('E51').Value = '  -0.37%  '
